# Recreated mantel correlograms with Euclidean distances.
# Update the "Mantel r" (column 3) and "p" (column 4) values in the
# correlogram table for each distance-class row (rows 2-12).
#
# Each table cell's Range.Text ends with a paragraph mark + cell mark
# (2 trailing characters), so we compute a sub-range that covers just
# the visible text and replace that - this keeps the edit strictly
# scoped to the target cell (Find/Replace on a cell Range was observed
# to occasionally jump to matches elsewhere in the document, which we
# want to avoid here).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $textLen = $rng.Text.Length - 2
    $textRng = $d.Range($rng.Start, $rng.Start + $textLen)
    $textRng.Text = $newText
}

# Row 2 - Distance class 1,500
Set-CellText $t 2 3 "0.017"
Set-CellText $t 2 4 "0.25"

# Row 3 - Distance class 4,500
Set-CellText $t 3 3 "-0.025"
Set-CellText $t 3 4 "0.308"

# Row 4 - Distance class 7,500
Set-CellText $t 4 3 "0.007"
Set-CellText $t 4 4 "0.5"

# Row 5 - Distance class 10,500
Set-CellText $t 5 3 "-0.014"
Set-CellText $t 5 4 "0.749"

# Row 6 - Distance class 13,500
Set-CellText $t 6 3 "-0.031"
Set-CellText $t 6 4 "0.624"

# Row 7 - Distance class 16,500
Set-CellText $t 7 3 "-0.033"
Set-CellText $t 7 4 "0.623"

# Row 8 - Distance class 19,500
Set-CellText $t 8 3 "-0.006"
Set-CellText $t 8 4 "0.999"

# Row 9 - Distance class 22,500 (only Mantel r changes; p stays "1")
Set-CellText $t 9 3 "-0.009"

# Row 10 - Distance class 25,500
Set-CellText $t 10 3 "0.036"
Set-CellText $t 10 4 "0.692"

# Row 11 - Distance class 28,500
Set-CellText $t 11 3 "0.023"
Set-CellText $t 11 4 "1"

# Row 12 - Distance class 31,500 (only Mantel r changes; p stays "1")
Set-CellText $t 12 3 "0.017"

Write-Output "edits applied"
